$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the current row 391, pushing the existing
# rows 391-409 down to 394-412 (matching the dimension growing from
# A1:R409 to A1:R412).
$ws.Rows("391:393").Insert()

# --- New row 391 ---
$ws.Range("A391").Value = 10
$ws.Range("B391").Value = "Vega Modelo de Temuco"
$ws.Range("C391").Value = "La Araucanía"
$ws.Range("D391").Value = 44753
$ws.Range("E391").Value = 9
$ws.Range("F391").Value = 100114014
$ws.Range("G391").Value = "Betarraga"
$ws.Range("H391").Value = "Sin especificar"
$ws.Range("I391").Value = "Primera"
$ws.Range("J391").Value = 50
$ws.Range("K391").Value = 10000
$ws.Range("L391").Value = 10000
$ws.Range("M391").Value = 10000
$ws.Range("N391").Value = "$/docena de paquetes"
$ws.Range("O391").Value = "Provincia de Cautín"
$ws.Range("P391").Value = 833
$ws.Range("Q391").Value = 12
$ws.Range("R391").Value = "Hortaliza"

# --- New row 392 ---
$ws.Range("A392").Value = 10
$ws.Range("B392").Value = "Vega Modelo de Temuco"
$ws.Range("C392").Value = "La Araucanía"
$ws.Range("D392").Value = 44753
$ws.Range("E392").Value = 9
$ws.Range("F392").Value = 100114014
$ws.Range("G392").Value = "Betarraga"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 100
$ws.Range("K392").Value = 8000
$ws.Range("L392").Value = 8000
$ws.Range("M392").Value = 8000
$ws.Range("N392").Value = "$/docena de paquetes"
$ws.Range("O392").Value = "Región del Maule"
$ws.Range("P392").Value = 667
$ws.Range("Q392").Value = 12
$ws.Range("R392").Value = "Hortaliza"

# --- New row 393 ---
$ws.Range("A393").Value = 10
$ws.Range("B393").Value = "Vega Modelo de Temuco"
$ws.Range("C393").Value = "La Araucanía"
$ws.Range("D393").Value = 44753
$ws.Range("E393").Value = 9
$ws.Range("F393").Value = 100114014
$ws.Range("G393").Value = "Betarraga"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 60
$ws.Range("K393").Value = 8000
$ws.Range("L393").Value = 9000
$ws.Range("M393").Value = 8500
$ws.Range("N393").Value = "$/saco 25 kilos"
$ws.Range("O393").Value = "Provincia de Cautín"
$ws.Range("P393").Value = 340
$ws.Range("Q393").Value = 25
$ws.Range("R393").Value = "Hortaliza"
